# "Updated IFR Lost comms"
# Slide 2 ("IFR Lost Comms") - shift the whole comms-lost/lost-comms-procedure
# diagram upward (the separator line under the title is removed, so everything
# below it moves up to close the gap).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# 1) Title textbox "IFR Lost Comms" moves up and slightly left.
$shp = $s.Shapes.Item("TextBox 2")
$shp.Left = 275.6023864746094
$shp.Top  = 75.34716796875

# 2) Remove the horizontal divider line under the title
#    ("Straight Connector 3", directly below the title textbox).
$s.Shapes.Item("Straight Connector 3").Delete()

# 3) "VFR / Cond." textbox moves up (x unchanged).
$shp = $s.Shapes.Item("TextBox 4")
$shp.Left = 294.4000244140625
$shp.Top  = 139.7100830078125

# 4) Arrow from "VFR / Cond." to "Continue VFR..." moves up (x unchanged).
$shp = $s.Shapes.Item("Straight Arrow Connector 5")
$shp.Left = 343.2631530761719
$shp.Top  = 161.82394409179688

# 5) Downward arrow from "VFR / Cond." moves up (x unchanged).
$shp = $s.Shapes.Item("Straight Arrow Connector 6")
$shp.Left = 318.83160400390625
$shp.Top  = 183.93780517578125

# 6) "Y" label textbox moves up (x unchanged).
$shp = $s.Shapes.Item("TextBox 7")
$shp.Left = 343.2631530761719
$shp.Top  = 140.85977172851562

# 7) "Continue VFR / Land a.s.a. practicable" textbox moves up (x unchanged).
$shp = $s.Shapes.Item("TextBox 10")
$shp.Left = 421.91064453125
$shp.Top  = 136.85008239746094

# 8) "Route: A.V.E.F" textbox moves up (x unchanged).
$shp = $s.Shapes.Item("TextBox 15")
$shp.Left = 280.4649658203125
$shp.Top  = 233.9844970703125

# 9) "Altitude: Max(...)" textbox moves up (x unchanged).
$shp = $s.Shapes.Item("TextBox 16")
$shp.Left = 280.4639587402344
$shp.Top  = 277.2103271484375

# 10) "N" label textbox moves up (x unchanged).
$shp = $s.Shapes.Item("TextBox 18")
$shp.Left = 317.58843994140625
$shp.Top  = 179.6507110595703

# 11) "Assigned > Vectors > Expected > Filed" textbox moves up (x unchanged).
$shp = $s.Shapes.Item("TextBox 8")
$shp.Left = 419.0035705566406
$shp.Top  = 223.50732421875
